$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$firstRow = 2
$lastRow = 91

# Snapshot the existing Date (col A) and Pages (col C) columns for every
# data row before overwriting anything.
$dates = @{}
$vals = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $dates[$r] = $ws.Cells.Item($r, 1).Value()
    $vals[$r] = $ws.Cells.Item($r, 3).Value()
}
# The GSC export rolled forward one day: the oldest day (2025-10-06) drops
# off the report and a new day (2026-01-04) is appended with a Pages count
# of 0.
$dates[$lastRow + 1] = "2026-01-04"
$vals[$lastRow + 1] = 0

# Stage column A as Text so the ISO-formatted date strings are written back
# as literal text (matching the source data) instead of Excel silently
# reinterpreting them as date serials.
$dateRange = $ws.Range("A$firstRow`:A$lastRow")
$dateRange.NumberFormat = "@"

# Every row shifts up to take on the next day's figures.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $dates[$r + 1]
    $ws.Cells.Item($r, 3).Value = $vals[$r + 1]
}

# Drop the temporary Text format so column A's style reverts to the sheet's
# default (unformatted) cell style, same as before the edit.
$dateRange.ClearFormats()
